# feat: add 2022-Q1 data
#
# - Insert a new worksheet "2022-Q1" right before "总计" (so tab order
#   becomes 2021-Q2, 2021-Q3, 2021-Q4, 2022-Q1, 总计).
# - Populate "2022-Q1" with the fund holdings table for that quarter.
# - Prepend a "2022-Q1" summary row to the "总计" sheet (shifting the
#   previously-existing rows down by one).

$wb = $excel.ActiveWorkbook

# Helper: write a value into a cell as TEXT (shared-string / inlineStr
# semantics) even when it looks numeric ("005112", "1.42", "0.0663", ...),
# without leaving behind a lingering NumberFormat/style on the cell -
# mirrors how the source workbook was produced (pandas/openpyxl, not by
# typing into Excel, which would otherwise coerce these to numbers and/or
# tag the cell with a Text number format).
function Set-TextValue {
    param($ws, $cellRef, $value)
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $ws.Range("ZZ999").Copy()
    $rng.PasteSpecial(-4122) | Out-Null
}

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet, placed immediately before "总计".
# ---------------------------------------------------------------------
# Duplicate "总计" itself (rather than Worksheets.Add()) so the new sheet
# starts out with the same sheetPr/pageMargins/etc. as its siblings, then
# wipe its data - it lands right after "2021-Q4", i.e. right before the
# original "总计".
$q4 = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Copy($null, $q4)
$newSheet = $wb.Worksheets.Item("总计 (2)")
$newSheet.Name = "2022-Q1"
$newSheet.Cells.Clear()

# Pull header/body formatting (font, border, alignment) from an existing
# quarter sheet that already uses the shared "s=2" style, so we reuse the
# same style index instead of minting new ones.
$styleDonor = $wb.Worksheets.Item("2021-Q3")

$styleDonor.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122) | Out-Null

$styleDonor.Range("A2").Copy()
$newSheet.Range("A2:A5").PasteSpecial(-4122) | Out-Null

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Row 2
$newSheet.Range("A2").Value = 0
Set-TextValue $newSheet "B2" "005112"
$newSheet.Range("C2").Value = "银华中证全指医药卫生指数增强"
Set-TextValue $newSheet "D2" "1.42"
Set-TextValue $newSheet "E2" "86.59"
Set-TextValue $newSheet "F2" "4.67"
Set-TextValue $newSheet "G2" "0.0663"
$newSheet.Range("H2").Value = 9

# Row 3
$newSheet.Range("A3").Value = 1
Set-TextValue $newSheet "B3" "004634"
$newSheet.Range("C3").Value = "新疆前海联合泳涛灵活配置混合A"
Set-TextValue $newSheet "D3" "1.33"
Set-TextValue $newSheet "E3" "89.65"
Set-TextValue $newSheet "F3" "4.25"
Set-TextValue $newSheet "G3" "0.0565"
$newSheet.Range("H3").Value = 8

# Row 4
$newSheet.Range("A4").Value = 2
Set-TextValue $newSheet "B4" "006235"
$newSheet.Range("C4").Value = "东方城镇消费主题混合"
Set-TextValue $newSheet "D4" "0.50"
Set-TextValue $newSheet "E4" "90.32"
Set-TextValue $newSheet "F4" "4.81"
Set-TextValue $newSheet "G4" "0.0240"
$newSheet.Range("H4").Value = 6

# Row 5
$newSheet.Range("A5").Value = 3
Set-TextValue $newSheet "B5" "007041"
$newSheet.Range("C5").Value = "新疆前海联合泳涛灵活配置混合C"
Set-TextValue $newSheet "D5" "0.00"
Set-TextValue $newSheet "E5" "89.65"
Set-TextValue $newSheet "F5" "4.25"
$newSheet.Range("G5").Value = 0
$newSheet.Range("H5").Value = 8

# ---------------------------------------------------------------------
# 2. Prepend a "2022-Q1" row to the "总计" sheet.
# ---------------------------------------------------------------------
# NOTE: worksheet variables in this host resolve by *position*, not a
# stable object handle - `$totalSheet` now points at whatever sheet sits
# at its old index (shifted by the inserts/copies above). Re-resolve
# "总计" by name before using it again.
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

# `Insert()` copies the row-above's formatting onto the new row (B2:D2
# picked up the header row's bold/border style) - strip that back to the
# default (no explicit style), matching the plain data rows below it.
$totalSheet.Range("ZZ999").Copy()
$totalSheet.Range("B2:D2").PasteSpecial(-4122) | Out-Null

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 4
$totalSheet.Range("D2").Value = 0.15

# Re-apply the "s=2" style to the new A2 cell (row insert does not carry
# over the donor row's A-column style automatically in all cases).
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122) | Out-Null
$totalSheet.Range("A2").Value = 0

# Renumber the remaining index column (A3:A6) to stay 1,2,3,...
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3

# ---------------------------------------------------------------------
# 3. Restore the original active sheet/selection (untouched by the diff).
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q2").Activate()
